$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the title string in A2 (shared string table entry)
$ws.Range("A2").Value = "Australia news anchors caught in Djokovic hot mic360p"

# Update the numeric values in row 2
$ws.Range("B2").Value = 62.06666666666667
$ws.Range("C2").Value = 10.1
$ws.Range("D2").Value = 0.9666666666666667
$ws.Range("E2").Value = 10.16666666666667
